$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue "D2" "66.803.45"
$ws.Range("E2").Value = "  +5.60%  "
Set-TextValue "D3" "3.531.17"
$ws.Range("E3").Value = "  +8.41%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue "D5" "566.05"
$ws.Range("E5").Value = "  +8.75%  "
Set-TextValue "D6" "187.02"
$ws.Range("E6").Value = "  +10.85%  "
$ws.Range("E7").Value = "  +10.44%  "
Set-TextValue "D8" "3.522.74"
$ws.Range("E8").Value = "  +8.27%  "
$ws.Range("E9").Value = "  +0.08%  "
Set-TextValue "D10" "0.646"
$ws.Range("E10").Value = "  +8.92%  "
Set-TextValue "D11" "0.158"
$ws.Range("E11").Value = "  +21.06%  "
Set-TextValue "D12" "56.02"
$ws.Range("E12").Value = "  +8.54%  "
Set-TextValue "D13" "0.0000282"
$ws.Range("E13").Value = "  +11.31%  "
Set-TextValue "D14" "9.48"
$ws.Range("E14").Value = "  +7.84%  "
Set-TextValue "D15" "4.089.42"
$ws.Range("E15").Value = "  +8.87%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D16" "18.84"
$ws.Range("E16").Value = "  +9.88%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "3.528.44"
$ws.Range("E17").Value = "  +8.78%  "
$ws.Range("E18").Value = "  +5.90%  "
Set-TextValue "D19" "66.864.34"
$ws.Range("E19").Value = "  +6.11%  "
Set-TextValue "D20" "12.19"
$ws.Range("E20").Value = "  +10.70%  "
$ws.Range("E21").Value = "  +7.04%  "
Set-TextValue "D22" "422.41"
$ws.Range("E22").Value = "  +13.15%  "
Set-TextValue "D23" "4.14"
$ws.Range("E23").Value = "  +14.22%  "
Set-TextValue "D24" "86.52"
$ws.Range("E24").Value = "  +7.54%  "
Set-TextValue "D25" "4.30"
$ws.Range("E25").Value = "  +4.74%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D26" "2.96"
$ws.Range("E26").Value = "  +11.42%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D27" "11.07"
$ws.Range("E27").Value = "  +0.61%  "
Set-TextValue "D28" "12.49"
$ws.Range("E28").Value = "  +13.02%  "
Set-TextValue "D29" "6.11"
$ws.Range("E29").Value = "  -0.27%  "
Set-TextValue "D30" "9.20"
$ws.Range("E30").Value = "  +15.62%  "
Set-TextValue "D31" "30.62"
$ws.Range("E31").Value = "  +8.63%  "
Set-TextValue "D32" "6.77"
$ws.Range("E32").Value = "  +4.76%  "
Set-TextValue "D33" "625.78"
$ws.Range("E33").Value = "  +1.52%  "
Set-TextValue "D34" "11.94"
$ws.Range("E34").Value = "  +8.22%  "
Set-TextValue "D35" "0.113"
$ws.Range("E35").Value = "  +9.38%  "
Set-TextValue "D36" "60.44"
$ws.Range("E36").Value = "  +7.85%  "
$ws.Range("E37").Value = "  +23.03%  "
Set-TextValue "D38" "0.0₃0827"
$ws.Range("E38").Value = "  +14.97%  "
Set-TextValue "D39" "38.57"
$ws.Range("E39").Value = "  +9.28%  "
$ws.Range("E40").Value = "  -0.02%  "
Set-TextValue "D41" "0.390"
$ws.Range("E41").Value = "  +5.31%  "
Set-TextValue "D42" "3.36"
$ws.Range("E42").Value = "  +9.06%  "
Set-TextValue "D43" "3.127.39"
$ws.Range("E43").Value = "  +11.08%  "
Set-TextValue "D44" "0.998"
$ws.Range("E44").Value = "  +0.03%  "
Set-TextValue "D45" "2.65"
$ws.Range("E45").Value = "  +2.70%  "
Set-TextValue "D46" "2.89"
$ws.Range("E46").Value = "  +13.41%  "
Set-TextValue "D47" "3.33"
$ws.Range("E47").Value = "  +12.66%  "
Set-TextValue "D48" "0.0422"
$ws.Range("E48").Value = "  +8.27%  "
Set-TextValue "D49" "2.74"
$ws.Range("E49").Value = "  +3.05%  "
$ws.Range("E50").Value = "  +8.91%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D51" "8.61"
$ws.Range("E51").Value = "  +11.09%  "
